$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: capitalize "unassigned" -> "Unassigned"
$ws.Range("B29").Value = "Unassigned"
$ws.Range("C29").Value = "Unassigned"
$ws.Range("D29").Value = "Unassigned"

# Row 41: capitalize "unassigned" -> "Unassigned"
$ws.Range("B41").Value = "Unassigned"
$ws.Range("C41").Value = "Unassigned"
$ws.Range("D41").Value = "Unassigned"

# Rows 42/43: swap ASV_IDs and update assignments (row 42 becomes Unassigned,
# row 43 becomes the former row 42's Homo sapiens/Human data)
$ws.Range("A42").Value = "c0a3f3ed23f04247d92740a9502f8b57"
$ws.Range("B42").Value = "Unassigned"
$ws.Range("C42").Value = "Unassigned"
$ws.Range("D42").Value = "Unassigned"

$ws.Range("A43").Value = "307c55294ffe3b8aa46fce358d55590e"
$ws.Range("B43").Value = "Homo sapiens"
$ws.Range("C43").Value = "Human"
$ws.Range("D43").Value = "Human"

# Rows 55/56: swap ASV_IDs and update assignments (row 55 becomes the former
# row 56's Centropristis striata/Black sea bass/Teleost Fish data, row 56 becomes Unassigned)
$ws.Range("A55").Value = "975b1dbdc7405f6e27bf63893e91e0ed"
$ws.Range("B55").Value = "Centropristis striata"
$ws.Range("C55").Value = "Black sea bass"
$ws.Range("D55").Value = "Teleost Fish"

$ws.Range("A56").Value = "5e733a21f67e541f28ed4bf4fe025044"
$ws.Range("B56").Value = "Unassigned"
$ws.Range("C56").Value = "Unassigned"
$ws.Range("D56").Value = "Unassigned"

# Row 60: capitalize "unassigned" -> "Unassigned"
$ws.Range("B60").Value = "Unassigned"
$ws.Range("C60").Value = "Unassigned"
$ws.Range("D60").Value = "Unassigned"
